$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: the fund name cell A20 changes from "科创债ETF平安" to "科创债ETF景顺"
# keeping the same fund code (B20 = 159400), but with a new bold font style.
$cell = $ws.Range("A20")
$cell.Value = "科创债ETF景顺"
$cell.Font.Bold = $true
$cell.Font.Name = "微软雅黑"
$cell.Font.Color = 2039583

# Restore the active selection to match the saved state.
$ws.Range("E6").Select()
